$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    if ($cell.Value -eq "5-11-2012-13") {
        $cell.Value = "2013-05-11"
    }
}
